# Apply "Add data for 2021-12-23" style update:
# - Rename sheet / title text from "through December 14" to "through December 15"
# - Update several carjacking counts across the grid (new data point added)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet tab name (was "Through 2021-12-14")
$ws.Name = "Through 2021-12-15"

# Update the header label in B1 (shared string text)
$ws.Range("B1").Value = "December 2021 (through December 15)"

# Cell value updates: column letter -> new value, keyed by row
$updates = @{
    "Z4"  = 2
    "AX4" = 3
    "Z6"  = 1
    "AX6" = 7
    "B7"  = 4
    "AL7" = 7
    "AX7" = 5
    "BJ7" = 4
    "BV7" = 2
    "B8"  = 5
    "N11" = 2
    "N12" = 2
    "B13" = 2
    "N13" = 1
    "B15" = 2
    "AX16" = 2
    "N18" = 3
    "AX18" = 1
    "B24" = 4
    "N24" = 5
    "B30" = 3
    "N34" = 1
    "B40" = 3
    "BJ43" = 2
    "B44" = 2
    "B54" = 2
    "N54" = 2
    "AL67" = 2
    "B70" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
